# Apply updated cryptocurrency price/volume figures to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.225.77"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.91%  "
$ws.Range("D3").Value = "'1.550.34"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -4.83%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'206.46"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.58%  "
$ws.Range("D6").Value = "'1.00"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "'0.475"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -5.68%  "
$ws.Range("E8").Value = "  -2.19%  "
$ws.Range("E9").Value = "  -3.85%  "
$ws.Range("D10").Value = "'17.76"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.84%  "
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("D12").Value = "'1.765.18"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.89%  "
$ws.Range("D13").Value = "'1.557.90"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.35%  "
$ws.Range("E14").Value = "  -5.25%  "
$ws.Range("E15").Value = "  -5.15%  "
$ws.Range("D16").Value = "'25.186.74"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.06%  "
$ws.Range("E17").Value = "  -4.40%  "
$ws.Range("D18").Value = "'58.54"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.58%  "
$ws.Range("D19").Value = "'1.00"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "'185.38"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.62%  "
$ws.Range("E21").Value = "  -3.56%  "
$ws.Range("E22").Value = "  -3.91%  "
$ws.Range("E23").Value = "  -4.30%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  -4.30%  "
$ws.Range("D26").Value = "'139.40"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.98%  "
$ws.Range("E27").Value = "  -5.14%  "
$ws.Range("E28").Value = "  -2.99%  "
$ws.Range("E30").Value = "  -6.92%  "
$ws.Range("E31").Value = "  -4.84%  "
$ws.Range("E32").Value = "  -4.07%  "
$ws.Range("E33").Value = "  -5.03%  "
$ws.Range("D34").Value = "'1.45"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.56%  "
$ws.Range("E35").Value = "  -4.05%  "
$ws.Range("D36").Value = "'1.084.59"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.93%  "
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").Value = "'0.493"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.64%  "
$ws.Range("E40").Value = "  -7.37%  "
$ws.Range("D41").Value = "'0.801"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.82%  "
$ws.Range("E42").Value = "  -11.42%  "
$ws.Range("D43").Value = "'92.53"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -5.76%  "
$ws.Range("D44").Value = "'5.04"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.68%  "
$ws.Range("D45").Value = "'1.680.98"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.84%  "
$ws.Range("D46").Value = "'0.0₆0112"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +13.46%  "
$ws.Range("D47").Value = "'52.21"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.09%  "
$ws.Range("E48").Value = "  -2.71%  "
$ws.Range("E49").Value = "  -5.89%  "
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("E51").Value = "  -2.02%  "
